# Recording of selenium execution: update Start time / End time / Time taken
# columns (E,F,G) for rows 2-6 on the active sheet with the latest run's
# timestamps, as captured by the selenium automation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "2023-07-12 10:50:01"
$ws.Range("F2").Value = "2023-07-12 10:51:33"
$ws.Range("G2").Value = "00:01:32"

$ws.Range("E3").Value = "2023-07-12 10:51:35"
$ws.Range("F3").Value = "2023-07-12 10:53:08"
$ws.Range("G3").Value = "00:01:33"

$ws.Range("E4").Value = "2023-07-12 10:53:09"
$ws.Range("F4").Value = "2023-07-12 10:54:41"
$ws.Range("G4").Value = "00:01:32"

$ws.Range("E5").Value = "2023-07-12 10:54:42"
$ws.Range("F5").Value = "2023-07-12 10:56:15"
$ws.Range("G5").Value = "00:01:33"

$ws.Range("E6").Value = "2023-07-12 10:56:20"
$ws.Range("F6").Value = "2023-07-12 11:05:31"
$ws.Range("G6").Value = "00:09:11"
